# Apply updated crypto price/volume data as scraped on Thu Sep  5 10:43:31 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'56.670.56"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.39%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.388.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.63%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  -0.06%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'504.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.74%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'132.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.57%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  +0.09%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.550"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.26%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'2.391.70"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.08%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.0974"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.69%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'  +0.55%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'  +1.38%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'4.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.29%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'2.810.66"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.55%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'56.596.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.60%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'21.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.62%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = "'  +1.36%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'2.359.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.19%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = "'  +0.09%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'  +0.41%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'308.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.03%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'  +1.22%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = "'  +0.19%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'5.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.91%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'65.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.50%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.13%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'0.382"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +3.13%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  +0.46%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "'  +2.86%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'176.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.36%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("E31").Value = "'  +2.07%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = "'  -0.30%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = "'  +1.16%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'5.85"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -4.07%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "'  +0.15%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  +0.29%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'17.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.22%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'  -2.24%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'3.81"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.80%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("B40").Value = "'SuiNetwork"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'0.821"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +6.17%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("B41").Value = "'OKB"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'36.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.36%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'  +1.10%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'130.81"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.92%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").Value = "'  +0.94%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'4.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.61%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.567"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.50%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'  +1.28%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'248.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.08%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = "'  -0.37%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.0210"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.59%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'17.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +7.25%  "
$ws.Range("E51").Style = "Normal"
